$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Second table (row 16): new entry, added before row 11 edits so that
# the "Added lava" shared string is interned before the others, matching
# the author's original edit order.
$ws.Range("Q16").Value2 = 44571
$ws.Range("Q16").NumberFormat = $ws.Range("Q15").NumberFormat
$ws.Range("R16").Value2 = "10.30 - 13.00"
$ws.Range("T16").Value2 = 2.5
$ws.Range("U16").Value2 = "Added lava"

# --- First table (row 11): new entry "After a (too) long break..."
$ws.Range("C11").Value2 = 44572
$ws.Range("C11").NumberFormat = $ws.Range("C10").NumberFormat
$ws.Range("D11").Value2 = "13.30 - 15.00"
$ws.Range("F11").Value2 = 1.5
$ws.Range("G11").Value2 = "After a (too) long break from this LO, trying to implement physics again"

# --- Remove the old intermediate subtotal formula in F14 (no longer needed
# now that row 11 holds real data and the subtotal moved to row 22).
$ws.Range("F14").ClearContents()

# --- Blank spacer row 17 (keeps date-style formatting like row 16 above it)
$ws.Range("Q17").NumberFormat = $ws.Range("Q15").NumberFormat

# --- Update the totals row (row 22) formulas to include the new rows
$ws.Range("F22").Formula = "=SUM(F3:F11)"
$ws.Range("T22").Formula = "=SUM(T3:T16)"

# --- Update the last selected cell to reflect where the user ended up
$ws.Range("F26").Select()
